$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2124183006535948
$ws.Cells.Item(2, 3).Value = 0.5359477124183006
$ws.Cells.Item(2, 10).Value = 0.02287581699346405
$ws.Cells.Item(2, 16).Value = 0.1339869281045752
$ws.Cells.Item(2, 19).Value = 0.09477124183006536

$ws.Cells.Item(3, 2).Value = 0.005714285714285714
$ws.Cells.Item(3, 3).Value = 0.04
$ws.Cells.Item(3, 10).Value = 0.02857142857142857
$ws.Cells.Item(3, 16).Value = 0.8
$ws.Cells.Item(3, 19).Value = 0.1257142857142857

$ws.Cells.Item(4, 10).Value = 0.02040816326530612
$ws.Cells.Item(4, 16).Value = 0.7346938775510204
$ws.Cells.Item(4, 19).Value = 0.2448979591836735

$ws.Cells.Item(6, 2).Value = 0.06422018348623854
$ws.Cells.Item(6, 4).Value = 0.009174311926605505
$ws.Cells.Item(6, 5).Value = 0.004587155963302753
$ws.Cells.Item(6, 6).Value = 0.04128440366972477
$ws.Cells.Item(6, 10).Value = 0.2752293577981652
$ws.Cells.Item(6, 15).Value = 0.02293577981651376
$ws.Cells.Item(6, 17).Value = 0.1100917431192661
$ws.Cells.Item(6, 18).Value = 0.1284403669724771
$ws.Cells.Item(6, 19).Value = 0.3440366972477064

$ws.Cells.Item(7, 2).Value = 0.08666666666666667
$ws.Cells.Item(7, 4).Value = 0.006666666666666667
$ws.Cells.Item(7, 5).Value = 0.006666666666666667
$ws.Cells.Item(7, 6).Value = 0.04
$ws.Cells.Item(7, 10).Value = 0.1333333333333333
$ws.Cells.Item(7, 15).Value = 0.01333333333333333
$ws.Cells.Item(7, 17).Value = 0.1533333333333333
$ws.Cells.Item(7, 18).Value = 0.1133333333333333
$ws.Cells.Item(7, 19).Value = 0.4466666666666667

$ws.Cells.Item(8, 2).Value = 0.06053268765133172
$ws.Cells.Item(8, 4).Value = 0.01452784503631961
$ws.Cells.Item(8, 6).Value = 0.04116222760290557
$ws.Cells.Item(8, 10).Value = 0.1428571428571428
$ws.Cells.Item(8, 15).Value = 0.02179176755447942
$ws.Cells.Item(8, 17).Value = 0.1694915254237288
$ws.Cells.Item(8, 18).Value = 0.12590799031477
$ws.Cells.Item(8, 19).Value = 0.423728813559322

$ws.Cells.Item(9, 2).Value = 0.05241935483870968
$ws.Cells.Item(9, 4).Value = 0.01209677419354839
$ws.Cells.Item(9, 6).Value = 0.03225806451612903
$ws.Cells.Item(9, 10).Value = 0.1451612903225807
$ws.Cells.Item(9, 15).Value = 0.02419354838709677
$ws.Cells.Item(9, 17).Value = 0.125
$ws.Cells.Item(9, 18).Value = 0.1411290322580645
$ws.Cells.Item(9, 19).Value = 0.4677419354838709

$ws.Cells.Item(10, 2).Value = 0.1081081081081081
$ws.Cells.Item(10, 4).Value = 0.02262727844123193
$ws.Cells.Item(10, 5).Value = 0.00251414204902577
$ws.Cells.Item(10, 6).Value = 0.06599622878692646
$ws.Cells.Item(10, 10).Value = 0.1194217473287241
$ws.Cells.Item(10, 15).Value = 0.01948460087994972
$ws.Cells.Item(10, 17).Value = 0.1961030798240101
$ws.Cells.Item(10, 18).Value = 0.1043368950345695
$ws.Cells.Item(10, 19).Value = 0.3614079195474544

$ws.Cells.Item(11, 7).Value = 0.1739130434782609
$ws.Cells.Item(11, 10).Value = 0.09881422924901186
$ws.Cells.Item(11, 11).Value = 0.225296442687747
$ws.Cells.Item(11, 12).Value = 0.5019762845849802

$ws.Cells.Item(12, 7).Value = 0.753968253968254
$ws.Cells.Item(12, 10).Value = 0.1984126984126984
$ws.Cells.Item(12, 11).Value = 0.01587301587301587
$ws.Cells.Item(12, 12).Value = 0.01587301587301587
$ws.Cells.Item(12, 19).Value = 0.01587301587301587

$ws.Cells.Item(13, 7).Value = 0.4827586206896552
$ws.Cells.Item(13, 10).Value = 0.4137931034482759
$ws.Cells.Item(13, 19).Value = 0.103448275862069

$ws.Cells.Item(15, 6).Value = 0.01526717557251908
$ws.Cells.Item(15, 8).Value = 0.1145038167938931
$ws.Cells.Item(15, 9).Value = 0.08015267175572519
$ws.Cells.Item(15, 10).Value = 0.4160305343511451
$ws.Cells.Item(15, 11).Value = 0.05725190839694656
$ws.Cells.Item(15, 15).Value = 0.07251908396946564
$ws.Cells.Item(15, 19).Value = 0.2442748091603053

$ws.Cells.Item(16, 6).Value = 0.02325581395348837
$ws.Cells.Item(16, 8).Value = 0.1395348837209302
$ws.Cells.Item(16, 9).Value = 0.08372093023255814
$ws.Cells.Item(16, 10).Value = 0.4790697674418605
$ws.Cells.Item(16, 11).Value = 0.06511627906976744
$ws.Cells.Item(16, 13).Value = 0.02325581395348837
$ws.Cells.Item(16, 14).Value = 0.004651162790697674
$ws.Cells.Item(16, 15).Value = 0.06046511627906977
$ws.Cells.Item(16, 19).Value = 0.1209302325581395

$ws.Cells.Item(17, 6).Value = 0.01754385964912281
$ws.Cells.Item(17, 8).Value = 0.1513157894736842
$ws.Cells.Item(17, 9).Value = 0.08991228070175439
$ws.Cells.Item(17, 10).Value = 0.4846491228070176
$ws.Cells.Item(17, 11).Value = 0.08333333333333333
$ws.Cells.Item(17, 13).Value = 0.01096491228070175
$ws.Cells.Item(17, 15).Value = 0.05263157894736842
$ws.Cells.Item(17, 19).Value = 0.1096491228070175

$ws.Cells.Item(18, 6).Value = 0.01677852348993289
$ws.Cells.Item(18, 8).Value = 0.1375838926174497
$ws.Cells.Item(18, 9).Value = 0.1174496644295302
$ws.Cells.Item(18, 10).Value = 0.4765100671140939
$ws.Cells.Item(18, 11).Value = 0.08053691275167785
$ws.Cells.Item(18, 13).Value = 0.01677852348993289
$ws.Cells.Item(18, 15).Value = 0.05704697986577181
$ws.Cells.Item(18, 19).Value = 0.09731543624161074

$ws.Cells.Item(19, 6).Value = 0.01630837657524092
$ws.Cells.Item(19, 8).Value = 0.1830985915492958
$ws.Cells.Item(19, 9).Value = 0.09710896960711639
$ws.Cells.Item(19, 10).Value = 0.4432913269088213
$ws.Cells.Item(19, 11).Value = 0.07412898443291327
$ws.Cells.Item(19, 13).Value = 0.01408450704225352
$ws.Cells.Item(19, 15).Value = 0.08080059303187546
$ws.Cells.Item(19, 19).Value = 0.09117865085248332
